$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds yearly rows 2008..2020 (rows 2-14).
# Target: drop the two oldest years (2008, 2009) and append a new
# year (2021) at the end, so the sheet ends up with 2010..2021 (rows 2-13).

# 1) Delete rows 2:3 (2008年, 2009年) - remaining rows shift up by two,
#    so former row 4 (2010年) becomes row 2, ... former row 14 (2020年)
#    becomes row 12.
$ws.Range("A2:S3").Delete()

# 2) Prime row 13 with the same cell layout/format as row 12 (copies the
#    "inlineStr" placeholders for blank columns E,F,G,J,P,S and the year
#    label style) before writing the new 2021年 figures into it.
$ws.Range("A12:S12").Copy()
$ws.Range("A13:S13").PasteSpecial(-4122)

# 3) Write the new 2021年 row of data.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 100128
$ws.Range("C13").Value = 310158
$ws.Range("D13").Value = 105606
$ws.Range("H13").Value = 2368.13
$ws.Range("I13").Value = 25536.95
$ws.Range("K13").Value = 5595.25
$ws.Range("L13").Value = 248363
$ws.Range("M13").Value = 1323304.5
$ws.Range("N13").Value = 123.21007
$ws.Range("O13").Value = 32524
$ws.Range("Q13").Value = 760999
$ws.Range("R13").Value = 270182
